$wb = $excel.ActiveWorkbook

# --- AddHeader sheet: insert a new "Division" column (with value 100) ---
# between "Asynchronous Processing" (B) and "Inventory Division" (C).
$wsHeader = $wb.Worksheets.Item("AddHeader")
$wsHeader.Activate()
$wsHeader.Columns.Item(3).Insert()
$wsHeader.Range("C1").Value = "Division"
$wsHeader.Range("C2").Value = 100
# Match the width of the "Asynchronous Processing" column (B) as closely as
# this engine's character-width -> XML-width rounding allows.
$wsHeader.Columns.Item(3).ColumnWidth = 21.584
$wsHeader.Range("C12").Select()

# --- CreditHold sheet: selection moved ---
$wsCreditHold = $wb.Worksheets.Item("CreditHold")
$wsCreditHold.Activate()
$wsCreditHold.Range("D34").Select()

# Restore the workbook's originally active tab (InventoryQuantity).
$wsActive = $wb.Worksheets.Item("InventoryQuantity")
$wsActive.Activate()
